$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the raw payroll period: rows 10-328 in column E should all read "04.21 - 05.04"
$ws.Range("E10:E328").Value = "04.21 - 05.04"
